# Saldo.xlsx update — refresh the "Export" balance list.
#
# Net effect (verified against the target diff):
#   - Remove PEDRO (004460487/100000), BRUNO (004452912/58052.5),
#     the old GUILHERME row (004948033/50360.8) and CAIO (004972351/38412.77).
#   - The CAIO row that used to hold 27029.32 (004512434) is replaced by a
#     LEDA row (002636063) with a refreshed balance of 22556.06.
#   - Remove the old LEDA row (002636063/13295.93) further down (duplicate
#     account, now consolidated above).
#   - Re-insert GUILHERME (004948033) with a refreshed balance of 9000.
#   - Insert a new BERNARDO row (005262440/3000) before ALESSANDRO.
#   - Remove the BLUEMETRIX row (001761119/2000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the four obsolete rows in one shot (old PEDRO, BRUNO, GUILHERME, CAIO).
$ws.Range("A3:A6").EntireRow.Delete()

# 2) The row that held CAIO/27029.32 has shifted up to row 4 — turn it into
#    the refreshed LEDA row.
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "002636063"
$ws.Cells.Item(4, 2).Value = "LEDA"
$ws.Cells.Item(4, 3).Value = 22556.06

# 3) The stale LEDA/13295.93 row has shifted up to row 6 — remove it.
$ws.Range("A6").EntireRow.Delete()

# 4) Re-insert GUILHERME with the refreshed balance at row 7.
$ws.Rows(7).Insert()
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "004948033"
$ws.Cells.Item(7, 2).Value = "GUILHERME"
$ws.Cells.Item(7, 3).Value = 9000

# 5) Insert the new BERNARDO row just before ALESSANDRO (row 11).
$ws.Rows(11).Insert()
$ws.Cells.Item(11, 1).NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "005262440"
$ws.Cells.Item(11, 2).Value = "BERNARDO"
$ws.Cells.Item(11, 3).Value = 3000

# 6) Drop the BLUEMETRIX row, now sitting at row 13.
$ws.Range("A13").EntireRow.Delete()

Write-Output "done"
